$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Category" header to "Vendor"
$ws.Range("B1").Value = "Vendor"

# Remove the trailing "Vendor" column (M1), which is now redundant
$ws.Range("M1").ClearContents()
